# Applies the "gh-pages build" restructuring:
#   - "2. Schedule" loses its last five sub-sessions, which become a new
#     top-level "3. " heading section (sessions 9-13 renumbered 2.6-2.10 -> 3.1-3.5).
#   - A new bookmark named "section" is introduced around the new heading
#     through the end of the document.

$d = $word.ActiveDocument

# --- 1. Insert the new "3. " Heading 1 paragraph right after
#        "2.5 Session 8: 19.06.25 - Talks 5 - 6" and before the
#        "session-9" bookmark / heading.
$session8Para = $d.Paragraphs(11)
$session8Para.Range.InsertParagraphAfter()
$newHeading = $d.Paragraphs(12)
$newHeading.Style = "Heading 1"
$newHeading.Range.Text = "3. "

# --- 2. Renumber the headings that move from section 2 to section 3.
$d.Content.Find.Execute(
    "2.6 Session 9: 26.06.25 - Talks 7 - 8", $true, $false, $false, $false, $false,
    $true, 1, $false, "3.1 Session 9: 26.06.25 - Talks 7 - 8", 2) | Out-Null

$d.Content.Find.Execute(
    "2.7 Session 10: 03.07.25 - Talks 9 - 10", $true, $false, $false, $false, $false,
    $true, 1, $false, "3.2 Session 10: 03.07.25 - Talks 9 - 10", 2) | Out-Null

$d.Content.Find.Execute(
    "2.8 Session 11: 10.07.25 - Talks 11 - 12", $true, $false, $false, $false, $false,
    $true, 1, $false, "3.3 Session 11: 10.07.25 - Talks 11 - 12", 2) | Out-Null

$d.Content.Find.Execute(
    "2.9 Session 12: 17.07.25 - Talks 13 - 14", $true, $false, $false, $false, $false,
    $true, 1, $false, "3.4 Session 12: 17.07.25 - Talks 13 - 14", 2) | Out-Null

$d.Content.Find.Execute(
    "2.10 Session 13: 24.07.25 - Talks 15 - 16", $true, $false, $false, $false, $false,
    $true, 1, $false, "3.5 Session 13: 24.07.25 - Talks 15 - 16", 2) | Out-Null

# --- 3. Wrap the new heading and everything after it (through the end of
#        the document, i.e. the rest of the old "schedule" bookmark plus
#        "references"/"refs") in a new bookmark named "section".
$sectionRange = $d.Range($session8Para.Range.End, $d.Content.End)
$d.Bookmarks.Add("section", $sectionRange)

Write-Output "Applied section split: '2.6-2.10' -> new '3.' heading with '3.1-3.5' sub-sessions."
